$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns H:L (1-indexed columns 8..12)
$headers = @("gray1", "gray2", "gray3", "gray4", "gray5")

# New grayscale image values for rows 2..6, columns H:L.
# These mirror the existing choice1..choice5 (columns B:F) image file names
# for each row, but pointing at the "_gray" variant of each image.
$grayData = @{
    2 = @("images/kiwi_gray.png",   "images/orange_gray.png", "images/grape_gray.png",  "images/banana_gray.png", "images/apple_gray.png")
    3 = @("images/kiwi_gray.png",   "images/apple_gray.png",  "images/grape_gray.png",  "images/banana_gray.png", "images/orange_gray.png")
    4 = @("images/apple_gray.png",  "images/banana_gray.png", "images/orange_gray.png", "images/kiwi_gray.png",   "images/grape_gray.png")
    5 = @("images/orange_gray.png", "images/apple_gray.png",  "images/grape_gray.png",  "images/banana_gray.png", "images/kiwi_gray.png")
    6 = @("images/banana_gray.png", "images/apple_gray.png",  "images/kiwi_gray.png",   "images/grape_gray.png",  "images/orange_gray.png")
}

# Write the header row (row 1), columns H..L -> column indices 8..12.
# New cells added to row 1 automatically pick up the bold/center header
# style already applied at the row level.
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 8 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Write the grayscale data rows 2..6, columns H..L -> column indices 8..12.
foreach ($row in 2..6) {
    $values = $grayData[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 8 + $i
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $values[$i]
        # Match the left-aligned style ("s=1") used by the sibling data cells.
        $cell.HorizontalAlignment = -4131
    }
}

# Move the active selection to I7, matching where the author's cursor
# ended up after entering the new data.
$ws.Range("I7").Select()
